# 10/08/2025 completing settlement details part2
#
# The merge-field placeholder "{{ cesantía_amount }}" lives in the
# "Cesantia" row of the additional-info table. Its inner text
# "cesantía_amount" is a single run. This change:
#   1. splits that run into three runs ("cesan" | "tia" | "_amount"),
#   2. drops the accent on the middle piece (tía -> tia) while leaving
#      the surrounding text untouched,
# all three runs keep identical run formatting.

$d = $word.ActiveDocument

# Locate the "cesantía_amount" run inside the merge-field braces.
$found = $d.Content
$ok = $found.Find.Execute("cesantía_amount", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $ok) {
    throw "could not find 'cesantía_amount' in the document"
}
$start = $found.Start
$end = $found.End

# "cesan" | "tía" (4th char removed of accent below) | "_amount"
$midStart = $start + 5
$midEnd = $start + 8

$mid = $d.Range($midStart, $midEnd)
$originalSize = $mid.Font.Size

# Applying a genuinely different direct-character-format value forces
# Word to split the run at the selection boundaries (it cannot keep
# reusing the single run once the properties on part of it diverge).
$mid.Font.Size = $originalSize + 1

# Replace the accented text with its plain-ASCII spelling now, while the
# run is still split out from its neighbours.
$mid.Text = "tia"

# Restore the original size so the three resulting runs share identical
# formatting again (Word keeps them as separate runs even though the
# properties now match their neighbours).
$mid.Font.Size = $originalSize
